$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item(1)

# --- Add AdminUsersPage sheet right after LoginPage ---
$wsAdmin = $wb.Worksheets.Add($null, $wsLogin)
$wsAdmin.Name = "AdminUsersPage"
$wsAdmin.Range("A1").Value = "UserName"
$wsAdmin.Range("B1").Value = "Password"
$wsAdmin.Range("A2").Value = "Rekha"
$wsAdmin.Range("B2").Value = "rekha"
$wsAdmin.Range("I21").Select()

# --- Add NewsPage sheet right after AdminUsersPage (becomes last/active sheet) ---
$wsNews = $wb.Worksheets.Add($null, $wsAdmin)
$wsNews.Name = "NewsPage"
$wsNews.Range("A1").Value = "NewsTitle"
$wsNews.Range("A2").Value = "Flash Sale Alert Get 50 percent off on select groceries this weekend only"

# Make NewsPage the active sheet/tab
$wsNews.Activate()
